# Corrected how event codes were shifted (errors in X1b)
#
# The data row had drifted one column to the right (B2 was blank, and the
# real values lived in C2:F2 with C2 holding a stale "2"). Remove the three
# now-unused trailing columns (D:F) so C2:F2 shift left into B2:C2, then
# write the corrected values: B2 = 2 (event code), C2 = 5 (corrected code,
# was erroneously 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the extra columns D:F - their contents (3, 4, 6) were shifted
# artifacts of the earlier mistake and are no longer needed; deleting them
# shifts C2:F2 left by one column.
$ws.Range("D1:F1").EntireColumn.Delete()

# Write the corrected event codes for row 2.
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 5

# Restore the expected active selection.
$ws.Range("D7").Select() | Out-Null
